# Adds two new measurement columns (LONG_MAN / ANCH_MAN) to the sheet,
# fixes a data-entry value in B12, upper-cases the "identificador de hoja"
# legend text, and appends two new legend rows describing the new columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New header cells E1 / F1
#    E1 mirrors the "ID/LONGITUD/ANCHO/AREA" header formatting (centered).
#    F1 mirrors the plain "font, no special alignment" formatting used by
#    the legend cells in row 29.
# ---------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "LONG_MAN"

$ws.Range("A29").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "ANCH_MAN"

# ---------------------------------------------------------------------
# 2. Data rows 2-26: new LONG_MAN (E) / ANCH_MAN (F) manual measurements.
#    Most values are dates (same "d.m" numeric format as LONGITUD/ANCHO);
#    a handful of rows hold plain numbers instead (no date format),
#    matching the exceptions already present in columns B/C.
# ---------------------------------------------------------------------
$longManValues = @{
  2 = 45885.0;  3 = 45820.0;  4 = 14.0;     5 = 13.0;     6 = 45912.0;
  7 = 45700.0;  8 = 45788.0;  9 = 45882.0;  10 = 45851.0; 11 = 45788.0;
  12 = 45760.0; 13 = 45850.0; 14 = 13.0;    15 = 45789.0; 16 = 45703.0;
  17 = 45703.0; 18 = 45732.0; 19 = 14.0;    20 = 13.0;    21 = 45880.0;
  22 = 45790.0; 23 = 45789.0; 24 = 45819.0; 25 = 45728.0; 26 = 45699.0
}
$longManPlain = @(4, 5, 14, 19, 20)

$anchManValues = @{
  2 = 45700.0;  3 = 45696.0;  4 = 45698.0;  5 = 45785.0;  6 = 8.0;
  7 = 45784.0;  8 = 45784.0;  9 = 45845.0;  10 = 45785.0; 11 = 45877.0;
  12 = 10.0;    13 = 45908.0; 14 = 10.0;    15 = 45907.0; 16 = 45786.0;
  17 = 45818.0; 18 = 45819.0; 19 = 45696.0; 20 = 45786.0; 21 = 45845.0;
  22 = 9.0;     23 = 8.0;     24 = 45784.0; 25 = 45845.0; 26 = 45814.0
}
$anchManPlain = @(6, 12, 14, 22, 23)

for ($row = 2; $row -le 26; $row++) {

  # --- E column (LONG_MAN) ---
  $eCell = $ws.Range("E$row")
  if ($longManPlain -contains $row) {
    # Plain-number look-alike: copy the "ID" column formatting (centered,
    # no date format).
    $ws.Range("A$row").Copy()
    $eCell.PasteSpecial(-4122)
  } else {
    # Date look-alike: copy the "LONGITUD" column formatting.
    $ws.Range("B$row").Copy()
    $eCell.PasteSpecial(-4122)
  }
  $eCell.Value = $longManValues[$row]

  # --- F column (ANCH_MAN) ---
  $fCell = $ws.Range("F$row")
  if ($anchManPlain -contains $row) {
    # Plain-number look-alike: same font/alignment as the row-29 legend,
    # no date format.
    $ws.Range("A29").Copy()
    $fCell.PasteSpecial(-4122)
  } else {
    # Date look-alike: same font/alignment as the row-29 legend, but with
    # the "d.m" date number format layered on top.
    $ws.Range("A29").Copy()
    $fCell.PasteSpecial(-4122)
    $fCell.NumberFormat = "d.m"
  }
  $fCell.Value = $anchManValues[$row]
}

# ---------------------------------------------------------------------
# 3. Data correction: B12 was mis-keyed as 23, should be 13.
# ---------------------------------------------------------------------
$ws.Range("B12").Value = 13.0

# ---------------------------------------------------------------------
# 4. Legend text: upper-case the sheet-identifier description, and add
#    two new legend rows describing the new columns (copy row-29 format).
# ---------------------------------------------------------------------
$ws.Range("B29").Value = "IDENTIFICADOR DE HOJA"

$ws.Range("A29:B29").Copy()
$ws.Range("A30:B31").PasteSpecial(-4122)

$ws.Range("A30").Value = "LONG_MAN"
$ws.Range("B30").Value = "MEDIDA MANUAL DE LONGITUD"
$ws.Range("A31").Value = "ANCH_MAN"
$ws.Range("B31").Value = "MEDIDA MANUAL DEL ANCHO"
